$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25.911653518676758
$ws.Range("C2").Value = 6.264367580413818
$ws.Range("D2").Value = 18.899436950683594
$ws.Range("E2").Value = 57.85714340209961
